$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Merge the "Assignment: " / "Calculator UI Design Wireframe" runs
#    into a single run holding the combined text.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Assignment: Calculator UI Design Wireframe", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Assignment: Calculator UI Design Wireframe", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Design 2 bullet: "...TextField is 45sp..." -> "...TextField is 60sp..."
#    Typed as two edits (change the "4"/"5" digits to "6"/"0"), which is
#    why the sentence ends up split across three runs afterwards.
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(
    "In Design 2 Text size for TextField is 45sp", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$base2 = $rng2.Start
$boundaryA2 = $base2 + 9     # right after "In Design"
$boundaryB2 = $base2 + 41    # right after the (soon to be) "60", before "sp"

$digitsRng2 = $d.Range($boundaryB2 - 2, $boundaryB2)
$digitsRng2.Text = "60"

# Drop transient bookmarks at the two edit boundaries (rightmost first) so
# the run split they force survives even after the markers are removed.
$markB2 = $d.Range($boundaryB2, $boundaryB2)
$d.Bookmarks.Add("TempMarkD2b", $markB2) | Out-Null
$markA2 = $d.Range($boundaryA2, $boundaryA2)
$d.Bookmarks.Add("TempMarkD2a", $markA2) | Out-Null
$d.Bookmarks.Item("TempMarkD2a").Delete() | Out-Null
$d.Bookmarks.Item("TempMarkD2b").Delete() | Out-Null

# ---------------------------------------------------------------------
# 3) Design 1 bullet: "...TextField is 45sp..." -> "...TextField is 65sp..."
#    The author placed the cursor right after the "4", deleted it and
#    typed "6" - leaving the insertion point (and so the document's
#    "_GoBack" bookmark) sitting between the new "6" and the old "5sp".
#    Word keeps only one "_GoBack" bookmark, so adding this one removes
#    the old one that used to sit near the end of the document.
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute(
    "In Design 1 Text size for TextField is 45sp", $true, $false, $false,
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$base1 = $rng1.Start
$boundaryText1 = $base1 + 26   # right before "TextField"
$digitPos1 = $base1 + 39       # position of the "4" in "45sp"

$digitRng1 = $d.Range($digitPos1, $digitPos1 + 1)
$digitRng1.Text = "6"

$caretPos1 = $digitPos1 + 1
$caretRng1 = $d.Range($caretPos1, $caretPos1)
$d.Bookmarks.Add("_GoBack", $caretRng1) | Out-Null

$tempRng1 = $d.Range($boundaryText1, $boundaryText1)
$d.Bookmarks.Add("TempMarkD1", $tempRng1) | Out-Null
$d.Bookmarks.Item("TempMarkD1").Delete() | Out-Null
